{"js": "// Bubble sort als H\u00dc\n// Insert a block of new \"Listenabsatz\" paragraphs right after the empty\n// paragraph that follows \"Sum, Power\" (and right before the run of\n// trailing empty paragraphs at the end of the document).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph: the first empty paragraph that comes right\n// after the paragraph whose text is \"Sum, Power\".\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  if (paragraphs.items[i].text === \"Sum, Power\" && paragraphs.items[i + 1].text === \"\") {\n    anchor = paragraphs.items[i + 1];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Anchor paragraph (empty paragraph after 'Sum, Power') not found\");\n}\n\n// New paragraph texts, in document order, each styled \"Listenabsatz\".\nconst newParagraphTexts = [\n  \"Prozess bekommt speicher bereich im Ram dort ist auch der Stack der Heap befindet sich irgendwo im Speicher\",\n  \"\",\n  \"Code Meta Daten Stack Speicher\",\n  \"Stack ist begrenzt kann also ausgehen\",\n  \"\",\n  \"Rekursion ist langsamer wie Iterativ\",\n  \"Wenn der Compiler dies umwandelt dann wird nicht die effektivste Variante genommen\",\n];\n\nlet insertAfter = anchor;\nfor (const text of newParagraphTexts) {\n  const p = insertAfter.insertParagraph(text, \"After\");\n  p.style = \"Listenabsatz\";\n  insertAfter = p;\n}\n\nawait context.sync();\n", "ps1": "# Bubble sort als H\u00dc\n# Insert a block of new \"Listenabsatz\" paragraphs right after the empty\n# paragraph that follows \"Sum, Power\" (and right before the run of\n# trailing empty paragraphs at the end of the document).\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph: the first empty paragraph that comes right\n# after the paragraph whose text is \"Sum, Power\".\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs($i).Range.Text.TrimEnd(\"`r\", \"`a\")\n  if ($t -eq \"Sum, Power\" -and $i -lt $count) {\n    $nextText = $d.Paragraphs($i + 1).Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($nextText -eq \"\") {\n      $anchorIndex = $i + 1\n      break\n    }\n  }\n}\n\nif ($anchorIndex -eq -1) {\n  throw \"Anchor paragraph (empty paragraph after 'Sum, Power') not found\"\n}\n\n# New paragraph texts, in document order, each styled \"Listenabsatz\".\n$newParagraphTexts = @(\n  \"Prozess bekommt speicher bereich im Ram dort ist auch der Stack der Heap befindet sich irgendwo im Speicher\",\n  \"\",\n  \"Code Meta Daten Stack Speicher\",\n  \"Stack ist begrenzt kann also ausgehen\",\n  \"\",\n  \"Rekursion ist langsamer wie Iterativ\",\n  \"Wenn der Compiler dies umwandelt dann wird nicht die effektivste Variante genommen\"\n)\n\n$insertAfterIndex = $anchorIndex\nforeach ($txt in $newParagraphTexts) {\n  $rng = $d.Paragraphs($insertAfterIndex).Range\n  $rng.InsertParagraphAfter()\n  $insertAfterIndex = $insertAfterIndex + 1\n  $newPara = $d.Paragraphs($insertAfterIndex)\n  $newPara.Style = \"Listenabsatz\"\n  if ($txt -ne \"\") {\n    $newPara.Range.Text = $txt\n  }\n}\n"}
